# ------------------------------------------------------------------
# combinations.xlsx edit: "All elements function so that we can get
# one of the figures for generalized pooling"
# ------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --------------------------------------------------------------
# 1) Tidy up leftover "no-fill" formatting on column G (rows 2-14)
#    and column M (several rows) that no longer carries any visible
#    formatting difference from the default style - clear it.
#    G2:G14 and M15 were empty (format-only) cells, so a full Clear
#    removes them outright; M10-M37 carry real values, so only their
#    formatting is cleared and the value is kept.
# --------------------------------------------------------------
$gRows = @(2,3,4,5,6,7,8,9,10,11,12,13,14)
foreach ($r in $gRows) {
    $ws.Range("G$r").Clear()
}

$ws.Range("M15").Clear()

$mRows = @(10,11,12,13,14,16,17,18,20,22,26,28,30,35,37)
foreach ($r in $mRows) {
    $ws.Range("M$r").ClearFormats()
}

# --------------------------------------------------------------
# 2) Fill in the "server" marker in column M for rows 108-115
#    (these experiments were run on the server, like many others
#    already marked earlier in the sheet).
# --------------------------------------------------------------
$serverRows = @(108,109,110,111,112,113,114,115)
foreach ($r in $serverRows) {
    $ws.Range("M$r").Value = "server"
}

# --------------------------------------------------------------
# 3) Add three placeholder index rows (147-149) - only the running
#    index in column J is populated, matching the blank rows already
#    present elsewhere in the index column.
# --------------------------------------------------------------
$ws.Range("J147").Value = 147
$ws.Range("J148").Value = 148
$ws.Range("J149").Value = 149

# --------------------------------------------------------------
# 4) Correct row 150 (Type: same, not small_compression) and
#    give it its running index too.
# --------------------------------------------------------------
$ws.Range("E150").Value = "same"
$ws.Range("J150").Value = 150

# --------------------------------------------------------------
# 5) Row 151: pooling is "gm" (generalized mean) not "mean", and
#    type is "same" not "small_compression"; add the Power value
#    and running index.
# --------------------------------------------------------------
$ws.Range("D151").Value = "gm"
$ws.Range("E151").Value = "same"
$ws.Range("G151").Value = 2
$ws.Range("J151").Value = 151

# --------------------------------------------------------------
# 6) New rows 152-156: continue the generalized-mean-pooling power
#    sweep (doubling G each row: 4, 8, 16, 32, 64) for the "same"
#    feature-number experiment, each with its running index.
# --------------------------------------------------------------
$newRows = @(
    @{Row=152; G=4},
    @{Row=153; G=8},
    @{Row=154; G=16},
    @{Row=155; G=32},
    @{Row=156; G=64}
)

foreach ($item in $newRows) {
    $r = $item.Row
    $ws.Range("A$r").Value = "BCE"
    $ws.Range("B$r").Value = 2
    $ws.Range("C$r").Value = "GCN"
    $ws.Range("D$r").Value = "gm"
    $ws.Range("E$r").Value = "same"
    $ws.Range("F$r").Value = 5
    $ws.Range("G$r").Value = $item.G
    $ws.Range("J$r").Value = $r
}

# --------------------------------------------------------------
# 7) Selection state: the author had scrolled down to the newly
#    extended area and selected the new Power column values.
# --------------------------------------------------------------
$ws.Activate()
$win = $wb.Windows.Item(1)
$win.ScrollRow = 136
$win.ScrollColumn = 1
$ws.Range("F150:F156").Select()
